$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 393, pushing existing rows 393-439 down to 394-440.
$ws.Rows.Item(393).Insert()

# Populate the newly inserted row 393 with the new weekly price record
# (same dimension/market/product columns as its neighbours; only the
# date + volume + price columns differ).
$ws.Cells.Item(393, 1).Value = 10
$ws.Cells.Item(393, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(393, 3).Value = "La Araucanía"
$ws.Cells.Item(393, 4).Value = 45212
$ws.Cells.Item(393, 5).Value = 9
$ws.Cells.Item(393, 6).Value = 100112052
$ws.Cells.Item(393, 7).Value = "Albahaca"
$ws.Cells.Item(393, 8).Value = "Sin especificar"
$ws.Cells.Item(393, 9).Value = "Primera"
$ws.Cells.Item(393, 10).Value = 30
$ws.Cells.Item(393, 11).Value = 5000
$ws.Cells.Item(393, 12).Value = 5000
$ws.Cells.Item(393, 13).Value = 5000
$ws.Cells.Item(393, 14).Value = "$/paquete"
$ws.Cells.Item(393, 15).Value = "Región Metropolitana"
$ws.Cells.Item(393, 16).Value = 5000
$ws.Cells.Item(393, 17).Value = 1
$ws.Cells.Item(393, 18).Value = "Hortaliza"
